# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Wed Apr 24 15:56:26 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.693.27'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.65%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.176.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.69%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.08%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.68%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.176.30'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.68%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.527'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.77%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.154'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.03%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.55'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.50%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.478'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.71%  '

# Row 13
$ws.Range('E13').Value = '  -6.47%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.33%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.699.31'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.61%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.765.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.73%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.178.53'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.65%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.15%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '482.53'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.21%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.48%  '

# Row 22
$ws.Range('E22').Value = '  -2.72%  '

# Row 23
$ws.Range('E23').Value = '  -3.99%  '

# Row 24
$ws.Range('E24').Value = '  -5.87%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.28%  '

# Row 26
$ws.Range('E26').Value = '  -0.20%  '

# Row 27
$ws.Range('E27').Value = '  -2.28%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.17%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.129'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +31.14%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.13%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.82%  '

# Row 32
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.73'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.59%  '

# Row 33
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.12%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.61%  '

# Row 35
$ws.Range('E35').Value = '  -6.21%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.08%  '

# Row 37
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.57%  '

# Row 38
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '54.43'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.97%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '474.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.42%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0734'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.71%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0404'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.91%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.125'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.54%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.52'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.09%  '

# Row 44
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.73%  '

# Row 45
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.900.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.86%  '

# Row 46
$ws.Range('E46').Value = '  -7.45%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.10%  '

# Row 49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.116'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.56%  '

# Row 50
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.68%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.08%  '
